$wb = $excel.ActiveWorkbook

# --- Quotes sheet: add two new quotes + Dutch translation for existing one ---
$quotes = $wb.Worksheets.Item("Quotes")

# Column widths: A narrower (author column), B much wider (holds long Dutch quote text)
$quotes.Columns.Item(1).ColumnWidth = 18.3
$quotes.Columns.Item(2).ColumnWidth = 140.1
$quotes.Columns.Item(2).BestFit = $false

# Row 3 (Marvin Minsky) gets its Dutch translation added in column B
$quotes.Range("B3").Value = "Zelden staan we stil bij het wonderbaarlijke feit dat we in ons leven nooit een 'echt serieuze fout maken, zoals een vork vol eten in je oog steken in plaats van in je mond, of een huis verlaten door een raam in plaats van een deur"

# Row 4: Allen Saunders quote (filled quote text first, then translation, source, author)
$quotes.Range("C4").Value = "Life is what happens to you, while you are busy making other plans"
$quotes.Range("B4").Value = "Het leven is wat je gebeurt, terwijl je andere plannen maakt"
$quotes.Range("D4").Value = "Acda en de Munnik"
$quotes.Range("A4").Value = "Allen Saunders"

# Row 5: Thomas Edison quote (quote text first, then author)
$quotes.Range("C5").Value = "But I have always found, when I was worrying, that the best thing to do was to put my mind upon something, work hard and forget what was troubling me. As a cure for worrying, work is better than whisky. Much better."
$quotes.Range("A5").Value = "Thomas Edison"

$quotes.Range("A6").Select()

# --- Make Planning the selected/active sheet ---
$planning = $wb.Worksheets.Item("Planning")
$planning.Select()
